# EPS 4.0 updates as of 3/15/2024
#
# 1. RACP sheet, cell B2: replace the formula
#      =ROUND(100/About!A11,0)   (cached 118)
#    with the hard-coded literal value 160.
# 2. About sheet: move the active selection from B15 to G17
#    (just a cursor-position / UI-state change).

$wb = $excel.ActiveWorkbook

$wsRacp  = $wb.Worksheets.Item("RACP")
$wsAbout = $wb.Worksheets.Item("About")

# Replace the formula in RACP!B2 with a plain literal value.
$wsRacp.Range("B2").Value = 160

# Update the remembered selection on the About sheet.
$wsAbout.Activate()
$wsAbout.Range("G17").Select() | Out-Null
